$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column F ("Viewers") is new: fill in the viewership figures for the
# final-predictions-vs-results dashboard (rows 2-33 of Sheet1).
$values = @{
    2 = 3560000
    3 = 582000
    4 = 466000
    5 = 0
    6 = 0
    7 = 0
    8 = 0
    9 = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 0
    21 = 0
    22 = 0
    23 = 0
    24 = 0
    25 = 0
    26 = 0
    27 = 0
    28 = 0
    29 = 0
    30 = 0
    31 = 0
    32 = 0
    33 = 0
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}

# Reproduce the sheet-view scroll/selection state from the edit: the
# bottom (split) pane is scrolled so row 21 is its first visible row,
# and the active cell moves from G38 to F34.
$win = $excel.ActiveWindow
$win.SplitColumn = 0
$win.SplitRow = 20
$win.FreezePanes = $false
$ws.Range("F34").Select()
